# Actualización automática 2025-06-26 15:05:09
#
# Updates sales figures for CHANDI ERAZO JOSUE after new sales were
# recorded for LAVABOS and PANELES PVC (client: CONSTANTE CAMACHO ARIANA
# ELIZABETH) in June.

$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO": per-client sales by product group ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("I2").Value = 157.15          # LAVABOS
$wsGrupo.Range("Q2").Value = 2322.22         # PANELES PVC

# Count of clients (out of 5) with sales in each group, row 7
$wsGrupo.Range("I7").Value = "1 de 5"        # LAVABOS now has a buyer
$wsGrupo.Range("Q7").Value = "1 de 5"        # PANELES PVC now has a buyer

# --- Sheet "VENTA MENSUAL": per-client monthly sales ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F2").Value = 2261.64       # junio, CONSTANTE CAMACHO ARIANA ELIZABETH
$wsMensual.Range("F7").Value = 3113.07       # junio, total

# --- Sheet "CUMPLIMIENTO MENSUAL": budget vs actual by product group ---
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# LAVABOS row (row 7)
$wsCumpl.Range("D7").Value = 157.15
$wsCumpl.Range("E7").Value = 342.85
$wsCumpl.Range("F7").Value = 0.3143

# PANELES PVC row (row 13)
$wsCumpl.Range("D13").Value = 2322.22
$wsCumpl.Range("E13").Value = -1839.22
$wsCumpl.Range("F13").Value = 4.807908902691511

# TOTAL row (row 17)
$wsCumpl.Range("D17").Value = 3113.07
$wsCumpl.Range("E17").Value = 10386.9318254209
$wsCumpl.Range("F17").Value = 0.2305977465971891
